$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E width: 18.6015625 -> 33.8359375
# (COM's ColumnWidth is quantized to whole pixels internally, so we dial in
# the input that rounds to the pixel bucket closest to the target width.)
$ws.Columns.Item(5).ColumnWidth = 33.142857142857144

# Row 37 restructuring: insert new D37/E37 ("Schuldenstandsquote" / "Ratio of government debt"),
# shifting the old D37/E37 ("Schuldenstand" / "Government debt") content into F37/G37
# (F37/G37 already held the same text, so the net effect is just D37/E37 getting new text).
$ws.Range("D37").Value = 'Schuldenstandsquote'
$ws.Range("E37").Value = 'Ratio of government debt'
$ws.Range("F37").Value = 'Schuldenstand'
$ws.Range("G37").Value = 'Government debt'

# Bulk text corrections across the sheet
$ws.Range("J2").Value = 'Anteil bis 2030 unter EU-27-Wert halten'
$ws.Range("K2").Value = 'keep the proportion below the EU-27-level by 2030'
$ws.Range("J3").Value = 'Anteil bis 2030 unter EU-27-Wert halten'
$ws.Range("K3").Value = 'keep the proportion below the EU-27-level by 2030'
$ws.Range("K4").Value = 'reduction to 70 kg/ha by 2030'
$ws.Range("J5").Value = 'Erhöhung des Anteils der landwirtschaftlichen Fläche auf 30 % bis 2030'
$ws.Range("K5").Value = 'increase the proportion of organically farmed agricultural land to 30% by 2030'
$ws.Range("D6").Value = 'Ernährungssicherung'
$ws.Range("E6").Value = 'Food Security'
$ws.Range("J6").Value = 'Steigerung des Anteils an den Gesamtausgaben bis 2030'
$ws.Range("K6").Value = 'shareas of total spending to be increased by 2030'
$ws.Range("K7").Value = 'to be reduced to 100 deaths per 100,000 inhabitants by 2030'
$ws.Range("K8").Value = 'to be reduced to 190 deaths per 100,000 inhabitants by 2030'
$ws.Range("K9").Value = 'reduction to 7% by 2030'
$ws.Range("K10").Value = 'reduction to 19% by 2030'
$ws.Range("K11").Value = 'increase to be permanently halted'
$ws.Range("K12").Value = 'increase to be permanently halted'
$ws.Range("K13").Value = 'reduction to 55% by 2030'
$ws.Range("K14").Value = 'adherence to the guideline value by 2030'
$ws.Range("K15").Value = 'expenditure to be increased by 2030'
$ws.Range("J16").Value = 'Verringerung des Anteils auf 9,5 % bis 2030'
$ws.Range("K16").Value = 'reduce the proportion to 9.5% by 2030'
$ws.Range("J17").Value = 'Steigerung des Anteils auf 55 % bis 2030'
$ws.Range("K17").Value = 'increase the proportion to 55% by 2030'
$ws.Range("K18").Value = 'increase to 35% by 2030'
$ws.Range("K19").Value = 'increase to 70% by 2030'
$ws.Range("J20").Value = 'Beibehaltung von 10 % bis 2030'
$ws.Range("K20").Value = 'maintaine 10% by 2030'
$ws.Range("J22").Value = 'gleichberechtigte Teilhabe bis 2025'
$ws.Range("K22").Value = 'equal-opportunity participation by 2025'
$ws.Range("D23").Value = 'Väterbeteiligung beim Elterngeld'
$ws.Range("K24").Value = 'to be increased by a third by 2030'
$ws.Range("J25").Value = 'Einhaltung oder Unterschreitung der Orientierungswerte bis 2030'
$ws.Range("K25").Value = 'not exceeding benchmark values by 2030'
$ws.Range("J26").Value = 'Einhaltung des Schwellenwertes von 50 mg/l bis 2030'
$ws.Range("K26").Value = 'compliance with the threshold value of 50 mg/l by 2030'
$ws.Range("K30").Value = 'increase by 2.1% per year'
$ws.Range("J31").Value = 'Senkung um 30 % bis 2030'
$ws.Range("K31").Value = 'reduction by 30% by 2030'
$ws.Range("J32").Value = 'Anstieg auf 30 % bis 2030'
$ws.Range("K32").Value = 'increase to 30% by 2030'
$ws.Range("D33").Value = 'Strom aus erneuerbaren Energiequellen'
$ws.Range("K33").Value = 'increase to 80% by 2030'
$ws.Range("D34").Value = 'Gesamtrohstoffproduktivität'
$ws.Range("E34").Value = 'Raw material input productivity'
$ws.Range("J34").Value = 'Beibehaltung des Trends der Jahre 2000 - 2010 bis 2030'
$ws.Range("K34").Value = 'trend of the years 2000-2010 to be maintained until 2030'
$ws.Range("J35").Value = 'unter 3 % des BIP'
$ws.Range("K35").Value = 'less than 3% of GDP'
$ws.Range("J36").Value = 'unter 0,5 % des BIP'
$ws.Range("K36").Value = 'less than 0.5% of GDP'
$ws.Range("J37").Value = 'max. 60 % des BIP'
$ws.Range("K37").Value = 'max. 60% of GDP'
$ws.Range("J38").Value = 'angemessene Entwicklung bis 2030'
$ws.Range("K38").Value = 'appropriate development by 2030'
$ws.Range("J39").Value = 'stetiges und angemessenes Wirtschaftswachstum'
$ws.Range("K39").Value = 'steady and appropriate economic growth'
$ws.Range("K40").Value = 'increase to 78% by 2030'
$ws.Range("K41").Value = 'increase to 60% by 2030'
$ws.Range("D42").Value = 'Mitglieder des Textilbündnisses'
$ws.Range("E42").Value = 'Members of the Textile Partnership'
$ws.Range("J42").Value = 'Steigerung der Anzahl bis 2030'
$ws.Range("K42").Value = 'increase number by 2030'
$ws.Range("J43").Value = 'jährlich mindestens 3,5 % des BIP bis 2025'
$ws.Range("K43").Value = 'at least 3.5% of GDP per year by 2025'
$ws.Range("J44").Value = 'flächendeckender Aufbau bis 2025'
$ws.Range("K44").Value = 'universal Roll-out by 2025'
$ws.Range("J45").Value = 'Erhöhung bei ausländischen Schulabsolvierenden und Angleichung an die Quote deutscher -absolvierenden bis 2030'
$ws.Range("K45").Value = 'increase among foreign school leavers and bring it into line with the rate of German school leavers by 2030'
$ws.Range("D46").Value = 'Gini-Koeffizient Einkommen nach Sozialtransfer'
$ws.Range("J46").Value = 'bis 2030 unterhalb des EU-27-Wertes'
$ws.Range("K46").Value = 'to be below the EU-27 figure by 2030'
$ws.Range("K47").Value = 'reduction to under 30 ha per day by 2030'
$ws.Range("K48").Value = 'reduce the loss'
$ws.Range("J49").Value = 'keine Verringerung'
$ws.Range("K49").Value = 'no reduction'
$ws.Range("K50").Value = 'reduction by 15-20% by 2030'
$ws.Range("K51").Value = 'reduction by 15-20% by 2030'
$ws.Range("D52").Value = 'Reisezeit mit öffentlichen Verkehrsmitteln'
$ws.Range("K52").Value = 'reduction'
$ws.Range("K53").Value = 'reduce to 13% by 2030'
$ws.Range("K54").Value = 'increase to 50 million by 2030'
$ws.Range("J55").Value = 'Steigerung des Marktanteils auf 34 % bis 2030'
$ws.Range("K55").Value = 'increase the market share to 34% by 2030'
$ws.Range("J56").Value = 'kontinuierliche Reduzierung'
$ws.Range("K56").Value = 'steady reduction'
$ws.Range("J57").Value = 'kontinuierliche Reduzierung'
$ws.Range("K57").Value = 'steady reduction'
$ws.Range("J58").Value = 'kontinuierliche Reduzierung'
$ws.Range("K58").Value = 'steady reduction'
$ws.Range("D59").Value = 'Umweltmanagement EMAS'
$ws.Range("E59").Value = 'EMAS eco-management'
$ws.Range("J59").Value = '5 000 Organisationsstandorte bis 2030'
$ws.Range("K59").Value = '5,000 locations of organisations by 2030'
$ws.Range("K60").Value = 'increase to 95% by 2020'
$ws.Range("J61").Value = 'signifikante Senkung'
$ws.Range("K61").Value = 'significantly reduce'
$ws.Range("K62").Value = 'reduce by 65% by 2030'
$ws.Range("K63").Value = 'increase to 6 bn euro by 2025'
$ws.Range("K64").Value = 'total nitrogen in the inflows below 2.6 mg/l'
$ws.Range("K65").Value = 'total nitrogen in the inflows below 2.8 mg/l'
$ws.Range("J66").Value = 'nachhaltige Bewirtschaftung nach dem MSY-Ansatz bis 2020'
$ws.Range("K66").Value = 'sustainable management in accordance with the MSY approachby 2020'
$ws.Range("K67").Value = 'reach the index value of 100 by 2030'
$ws.Range("K68").Value = 'reduction by 35% by 2030'
$ws.Range("D69").Value = 'a) Wälder unter REDD+-Regelwerk'
$ws.Range("E69").Value = 'a) Forests under REDD+ rulebook'
$ws.Range("K69").Value = 'increase payments by 2030'
$ws.Range("K70").Value = 'increase payments by 2030'
$ws.Range("K71").Value = 'reduce to less than 6,500 per 100,000 inhabitants by 2030'
$ws.Range("J72").Value = 'mindestens 15 Projekte pro Jahr bis 2030'
$ws.Range("K72").Value = 'at least 15 projects per year by 2030'
$ws.Range("K73").Value = 'improvement by 2030'
$ws.Range("K74").Value = 'improvement by 2030'
$ws.Range("J75").Value = 'Steigerung des Anteils auf 0,7 % des BNE bis 2030'
$ws.Range("K75").Value = 'increase to 0.7% of GNI by 2030'
$ws.Range("K76").Value = 'increase by 10% from 2015 to 2020, then stabilised'
$ws.Range("J77").Value = 'Steigerung des Anteils um 100 % bis 2030 gegenüber 2014'
$ws.Range("K77").Value = 'increase by 100 % by 2030 compared to 2014'
